$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain value updates (text and numeric columns) ---
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = "T.J. McConnell"
$ws.Range("D4").Value = "PG"
$ws.Range("E4").Value = "6-1"
$ws.Range("F4").Value = 190
$ws.Range("G4").Value = "March 25, 1992"
$ws.Range("J4").Value = "Duquesne, Arizona"
$ws.Range("K4").Value = "https://www.basketball-reference.com/players/m/mccontj01.html"
$ws.Range("B5").Value = 23
$ws.Range("C5").Value = "Aaron Nesmith"
$ws.Range("D5").Value = "SF"
$ws.Range("E5").Value = "6-5"
$ws.Range("F5").Value = 215
$ws.Range("G5").Value = "October 16, 1999"
$ws.Range("J5").Value = "Vanderbilt"
$ws.Range("K5").Value = "https://www.basketball-reference.com/players/n/nesmiaa01.html"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "Andrew Nembhard"
$ws.Range("D6").Value = "SG"
$ws.Range("F6").Value = 193
$ws.Range("G6").Value = "January 16, 2000"
$ws.Range("H6").Value = "ca"
$ws.Range("I6").Value = "R"
$ws.Range("J6").Value = "Florida, Gonzaga"
$ws.Range("K6").Value = "https://www.basketball-reference.com/players/n/nembhan01.html"
$ws.Range("B7").Value = 25
$ws.Range("C7").Value = "Jalen Smith"
$ws.Range("D7").Value = "PF"
$ws.Range("E7").Value = "6-10"
$ws.Range("F7").Value = 215
$ws.Range("G7").Value = "March 16, 2000"
$ws.Range("H7").Value = "us"
$ws.Range("J7").Value = "Maryland"
$ws.Range("K7").Value = "https://www.basketball-reference.com/players/s/smithja04.html"
$ws.Range("B8").Value = 33
$ws.Range("C8").Value = "Myles Turner"
$ws.Range("D8").Value = "C"
$ws.Range("E8").Value = "6-11"
$ws.Range("F8").Value = 250
$ws.Range("G8").Value = "March 24, 1996"
$ws.Range("H8").Value = "us"
$ws.Range("J8").Value = "Texas"
$ws.Range("K8").Value = "https://www.basketball-reference.com/players/t/turnemy01.html"
$ws.Range("B9").Value = 12
$ws.Range("C9").Value = "Oshae Brissett"
$ws.Range("E9").Value = "6-7"
$ws.Range("F9").Value = 210
$ws.Range("G9").Value = "June 20, 1998"
$ws.Range("H9").Value = "ca"
$ws.Range("J9").Value = "Syracuse"
$ws.Range("K9").Value = "https://www.basketball-reference.com/players/b/brissos01.html"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = "Tyrese Haliburton"
$ws.Range("D10").Value = "PG"
$ws.Range("E10").Value = "6-5"
$ws.Range("F10").Value = 185
$ws.Range("G10").Value = "February 29, 2000"
$ws.Range("J10").Value = "Iowa State"
$ws.Range("K10").Value = "https://www.basketball-reference.com/players/h/halibty01.html"
$ws.Range("B13").Value = 8
$ws.Range("C13").Value = "Trevelin Queen (TW)"
$ws.Range("E13").Value = "6-6"
$ws.Range("F13").Value = 190
$ws.Range("G13").Value = "February 25, 1997"
$ws.Range("J13").Value = "College of Marin, New Mexico Military Institute, New Mexico State"
$ws.Range("K13").Value = "https://www.basketball-reference.com/players/q/queentr01.html"
$ws.Range("B14").Value = 10
$ws.Range("C14").Value = "Kendall Brown (TW)"
$ws.Range("D14").Value = "SG"
$ws.Range("E14").Value = "6-8"
$ws.Range("F14").Value = 205
$ws.Range("G14").Value = "May 11, 2003"
$ws.Range("H14").Value = "us"
$ws.Range("I14").Value = "R"
$ws.Range("J14").Value = "Baylor"
$ws.Range("K14").Value = "https://www.basketball-reference.com/players/b/brownke03.html"
$ws.Range("B15").Value = 27
$ws.Range("C15").Value = "Daniel Theis"
$ws.Range("D15").Value = "C"
$ws.Range("E15").Value = "6-8"
$ws.Range("F15").Value = 245
$ws.Range("G15").Value = "April 4, 1992"
$ws.Range("H15").Value = "de"
$ws.Range("K15").Value = "https://www.basketball-reference.com/players/t/theisda01.html"
$ws.Range("C16").Value = "George Hill"
$ws.Range("D16").Value = "PG"
$ws.Range("E16").Value = "6-4"
$ws.Range("F16").Value = 188
$ws.Range("G16").Value = "May 4, 1986"
$ws.Range("J16").Value = "IUPUI"
$ws.Range("K16").Value = "https://www.basketball-reference.com/players/h/hillge01.html"
$ws.Range("C17").Value = "Serge Ibaka"
$ws.Range("D17").Value = "F-C"
$ws.Range("E17").Value = "6-10"
$ws.Range("F17").Value = 235
$ws.Range("G17").Value = "September 18, 1989"
$ws.Range("H17").Value = "cg"
$ws.Range("K17").Value = "https://www.basketball-reference.com/players/i/ibakase01.html"
$ws.Range("C18").Value = "Jordan Nwora"
$ws.Range("D18").Value = "SF"
$ws.Range("F18").Value = 225
$ws.Range("G18").Value = "September 9, 1998"
$ws.Range("H18").Value = "us"
$ws.Range("J18").Value = "Louisville"
$ws.Range("K18").Value = "https://www.basketball-reference.com/players/n/nworajo01.html"

# --- Clear cells that should become empty ---
$ws.Range("J15").Value = $null
$ws.Range("B16").Value = $null
$ws.Range("B17").Value = $null
$ws.Range("J17").Value = $null

# --- Numeric-looking strings in Exp column: must remain text (t="s") ---
# Use TEXT() formula then Copy/PasteSpecial(values) to bake a clean shared-string
# value without adding a NumberFormat style or leaving a formula behind.
$r = $ws.Range("I4")
$r.Formula = "=TEXT(7,""0"")"
$r.Copy()
$r.PasteSpecial(-4163)
$r = $ws.Range("I5")
$r.Formula = "=TEXT(2,""0"")"
$r.Copy()
$r.PasteSpecial(-4163)
$r = $ws.Range("I7")
$r.Formula = "=TEXT(2,""0"")"
$r.Copy()
$r.PasteSpecial(-4163)
$r = $ws.Range("I8")
$r.Formula = "=TEXT(7,""0"")"
$r.Copy()
$r.PasteSpecial(-4163)
$r = $ws.Range("I9")
$r.Formula = "=TEXT(3,""0"")"
$r.Copy()
$r.PasteSpecial(-4163)
$r = $ws.Range("I10")
$r.Formula = "=TEXT(2,""0"")"
$r.Copy()
$r.PasteSpecial(-4163)
$r = $ws.Range("I15")
$r.Formula = "=TEXT(5,""0"")"
$r.Copy()
$r.PasteSpecial(-4163)
$r = $ws.Range("I16")
$r.Formula = "=TEXT(14,""0"")"
$r.Copy()
$r.PasteSpecial(-4163)
$r = $ws.Range("I17")
$r.Formula = "=TEXT(13,""0"")"
$r.Copy()
$r.PasteSpecial(-4163)
$r = $ws.Range("I18")
$r.Formula = "=TEXT(2,""0"")"
$r.Copy()
$r.PasteSpecial(-4163)

$excel.CutCopyMode = 0
